$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the special header style (bold font, thin border, centered/top alignment)
# from the header row A1:C1 so it reverts back to the default "Normal" style.
$ws.Range("A1:C1").Style = "Normal"

# Delete the row for "ゴシキセイガイインコ。" (row 252). Excel shifts all the
# following rows up by one, so the old row 253 ("モモンガ。") becomes row 252, etc.
$ws.Rows(252).Delete()
